# Finished Week 13 logging
# Update row 3 (Road "R" totals) on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 185
$wsOff.Range("C3").Value = 119
$wsOff.Range("D3").Value = 40
$wsOff.Range("E3").Value = 20
$wsOff.Range("G3").Value = 3

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 205
$wsDef.Range("C3").Value = 155
$wsDef.Range("D3").Value = 35
$wsDef.Range("E3").Value = 14
